$d = $word.ActiveDocument

# --- Context -----------------------------------------------------------
# The document is a single paragraph containing one run of text ("Fuck").
# The target edit turns it into the (centered, Times New Roman, 26pt)
# title "SortOfSort Method Performance Analysis", split into two runs:
#   - "SortOfSort"                    (wrapped in a spellcheck proofErr
#                                       bracket, same as Word's proofer
#                                       would add for a word it doesn't
#                                       recognise)
#   - " Method Performance Analysis"
# while leaving the trailing "_GoBack" bookmark and the section
# properties untouched.

$para = $d.Paragraphs(1)

# Recover the paragraph's own identity attributes (w14:paraId / w14:textId
# / w:rsidR / w:rsidRDefault, ...) straight from the live package so the
# rewritten paragraph continues to look like an in-place edit of the same
# paragraph instead of a brand new one. Falls back to a bare <w:p> if for
# some reason none are present.
$currentXml = $d.Content.WordOpenXML
if ($currentXml -match '<w:p\s[^>]*>') {
    $pOpenTag = $matches[0]
} else {
    $pOpenTag = '<w:p>'
}

# Run/paragraph-mark formatting shared by every run plus the paragraph
# mark itself: Times New Roman, 26pt (w:sz/w:szCs are in half-points).
$rPr = '<w:rPr>' +
         '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
         '<w:sz w:val="52"/>' +
         '<w:szCs w:val="52"/>' +
       '</w:rPr>'

$titleFirst = 'SortOfSort'
$titleRest = ' Method Performance Analysis'

$newParaXml = $pOpenTag +
  '<w:pPr><w:jc w:val="center"/>' + $rPr + '</w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>' + $titleFirst + '</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $titleRest + '</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
             '<w:document ' +
               'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
               'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
               'xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" ' +
               'mc:Ignorable="w14">' +
               '<w:body>' + $newParaXml + '</w:body>' +
             '</w:document>' +
           '</pkg:xmlData>' +
         '</pkg:part>' +
       '</pkg:package>'

# InsertXML replaces the contents of the range it's invoked on. Calling it
# on $d.Content (rather than the whole package) swaps in the rebuilt
# paragraph while leaving the trailing section-properties markup (sectPr)
# - which lives outside Content - untouched.
[void]$d.Content.InsertXML($pkg)
